# MF attention new fit result
# Refit the "MF_attention" model after dropping the forget_attention parameter
# (model now has 5 free parameters instead of 6). Updates the per-participant
# fit statistics (nlld, alpha, tau, gamma, eta, alpha_attention), removes the
# forget_attention column, and fixes up the AICc formula's parameter count.
# The ranking-indicator flags (AA:AH) and the row-39 totals are formulas that
# recompute automatically once the dependent cells below change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MF_fit_result")

# --- Updated per-participant fit statistics for the refit MF_attention model ---
# Columns: CH=nlld(86) CI=alpha(87) CJ=tau(88) CK=gamma(89) CL=eta(90) CM=alpha_attention(91)

$ws.Cells.Item(3, 86).Value = 947.23047099186
$ws.Cells.Item(3, 87).Value = 0.0478949097860716
$ws.Cells.Item(3, 88).Value = 0.610605131173725
$ws.Cells.Item(3, 89).Value = 1
$ws.Cells.Item(3, 90).Value = 0.00210050007491859
$ws.Cells.Item(3, 91).Value = 0.00327048559491657

$ws.Cells.Item(4, 86).Value = 1633.36922855364
$ws.Cells.Item(4, 87).Value = 0.0129945970923526
$ws.Cells.Item(4, 88).Value = 0.206161956628918
$ws.Cells.Item(4, 89).Value = 1
$ws.Cells.Item(4, 90).Value = 0
$ws.Cells.Item(4, 91).Value = 0.11405463407296

$ws.Cells.Item(5, 86).Value = 1493.61951217595
$ws.Cells.Item(5, 87).Value = 0.0298128816389625
$ws.Cells.Item(5, 88).Value = 0.206521892285656
$ws.Cells.Item(5, 89).Value = 1
$ws.Cells.Item(5, 90).Value = 0.00182439930686735
$ws.Cells.Item(5, 91).Value = 1

$ws.Cells.Item(6, 86).Value = 640.804782690471
$ws.Cells.Item(6, 87).Value = 0.0249040032402655
$ws.Cells.Item(6, 88).Value = 0.816012113961024
$ws.Cells.Item(6, 89).Value = 1
$ws.Cells.Item(6, 90).Value = 0.05
$ws.Cells.Item(6, 91).Value = 0.0358735549701983

$ws.Cells.Item(7, 86).Value = 1923.4562765736
$ws.Cells.Item(7, 87).Value = 0.0527091915872304
$ws.Cells.Item(7, 88).Value = 0.0465101310121665
$ws.Cells.Item(7, 89).Value = 1
$ws.Cells.Item(7, 90).Value = 0.00122264857928721
$ws.Cells.Item(7, 91).Value = 0.73735174440477

$ws.Cells.Item(8, 86).Value = 784.95645056098
$ws.Cells.Item(8, 87).Value = 0.0251021495502693
$ws.Cells.Item(8, 88).Value = 0.316199311574376
$ws.Cells.Item(8, 89).Value = 1
$ws.Cells.Item(8, 90).Value = 0.00840450029503101
$ws.Cells.Item(8, 91).Value = 0.205311174887297

$ws.Cells.Item(9, 86).Value = 1048.37243882714
$ws.Cells.Item(9, 87).Value = 0.0158279136912241
$ws.Cells.Item(9, 88).Value = 0.228505306976791
$ws.Cells.Item(9, 89).Value = 1
$ws.Cells.Item(9, 90).Value = 0.00131140419075448
$ws.Cells.Item(9, 91).Value = 0.167821559112807

$ws.Cells.Item(10, 86).Value = 1207.09623983972
$ws.Cells.Item(10, 87).Value = 0.0481774771221039
$ws.Cells.Item(10, 88).Value = 0.311817453534213
$ws.Cells.Item(10, 89).Value = 1
$ws.Cells.Item(10, 90).Value = 0.00668339698135319
$ws.Cells.Item(10, 91).Value = 0.088639587226111

$ws.Cells.Item(11, 86).Value = 925.105794592838
$ws.Cells.Item(11, 87).Value = 0.0331190384991378
$ws.Cells.Item(11, 88).Value = 0.302534946827141
$ws.Cells.Item(11, 89).Value = 1
$ws.Cells.Item(11, 90).Value = 0.00152351904441756
$ws.Cells.Item(11, 91).Value = 0.769694108208045

$ws.Cells.Item(12, 86).Value = 383.529276539232
$ws.Cells.Item(12, 87).Value = 0.0364854449095222
$ws.Cells.Item(12, 88).Value = 0.414179973891508
$ws.Cells.Item(12, 89).Value = 1
$ws.Cells.Item(12, 90).Value = 0.00940574024828402
$ws.Cells.Item(12, 91).Value = 0.552856665679566

$ws.Cells.Item(13, 86).Value = 940.008462718011
$ws.Cells.Item(13, 87).Value = 0.0253884644975633
$ws.Cells.Item(13, 88).Value = 0.27742123835669
$ws.Cells.Item(13, 89).Value = 1
$ws.Cells.Item(13, 90).Value = 0.00592447108905228
$ws.Cells.Item(13, 91).Value = 0.877262169941599

$ws.Cells.Item(14, 86).Value = 843.681194636191
$ws.Cells.Item(14, 87).Value = 0.0242183318287607
$ws.Cells.Item(14, 88).Value = 0.402285636191517
$ws.Cells.Item(14, 89).Value = 1
$ws.Cells.Item(14, 90).Value = 0.002807595019862
$ws.Cells.Item(14, 91).Value = 0.0275777796249967

$ws.Cells.Item(15, 86).Value = 686.855473129504
$ws.Cells.Item(15, 87).Value = 0.0222593126047024
$ws.Cells.Item(15, 88).Value = 0.314602966973673
$ws.Cells.Item(15, 89).Value = 1
$ws.Cells.Item(15, 90).Value = 0.0059929900014293
$ws.Cells.Item(15, 91).Value = 0.177724172191575

$ws.Cells.Item(16, 86).Value = 904.869411381784
$ws.Cells.Item(16, 87).Value = 0.0245906928592451
$ws.Cells.Item(16, 88).Value = 0.392534787875928
$ws.Cells.Item(16, 89).Value = 1
$ws.Cells.Item(16, 90).Value = 0.029201555757164
$ws.Cells.Item(16, 91).Value = 0.0770110195501512

$ws.Cells.Item(17, 86).Value = 1591.48186549099
$ws.Cells.Item(17, 87).Value = 0.0399910879170805
$ws.Cells.Item(17, 88).Value = 0.305316982317317
$ws.Cells.Item(17, 89).Value = 1
$ws.Cells.Item(17, 90).Value = 0.00179466317580028
$ws.Cells.Item(17, 91).Value = 0.669695233765021

$ws.Cells.Item(18, 86).Value = 1329.45431144281
$ws.Cells.Item(18, 87).Value = 0.0449822022855921
$ws.Cells.Item(18, 88).Value = 0.278539978865647
$ws.Cells.Item(18, 89).Value = 1
$ws.Cells.Item(18, 90).Value = 0.0132083332266635
$ws.Cells.Item(18, 91).Value = 0.082212173270954

$ws.Cells.Item(19, 86).Value = 443.960912945967
$ws.Cells.Item(19, 87).Value = 0.0139700477798757
$ws.Cells.Item(19, 88).Value = 0.419394919907785
$ws.Cells.Item(19, 89).Value = 1
$ws.Cells.Item(19, 90).Value = 0.00290192378862043
$ws.Cells.Item(19, 91).Value = 0.676607064569957

$ws.Cells.Item(20, 86).Value = 1003.41066784377
$ws.Cells.Item(20, 87).Value = 0.0283711713409262
$ws.Cells.Item(20, 88).Value = 0.290063599275902
$ws.Cells.Item(20, 89).Value = 1
$ws.Cells.Item(20, 90).Value = 0.00853476731684828
$ws.Cells.Item(20, 91).Value = 0.0517140015772316

$ws.Cells.Item(21, 86).Value = 1892.55708582335
$ws.Cells.Item(21, 87).Value = 0.110070757269657
$ws.Cells.Item(21, 88).Value = 0.0530142085743497
$ws.Cells.Item(21, 89).Value = 1
$ws.Cells.Item(21, 90).Value = 0.0066763516156531
$ws.Cells.Item(21, 91).Value = 0.991217086011442

$ws.Cells.Item(22, 86).Value = 496.131125558313
$ws.Cells.Item(22, 87).Value = 0.0252313511638338
$ws.Cells.Item(22, 88).Value = 0.308036519162832
$ws.Cells.Item(22, 89).Value = 1
$ws.Cells.Item(22, 90).Value = 0.00790921303236013
$ws.Cells.Item(22, 91).Value = 0.859703671102239

$ws.Cells.Item(23, 86).Value = 1109.84363670359
$ws.Cells.Item(23, 87).Value = 0.0170443969934115
$ws.Cells.Item(23, 88).Value = 0.273011643923626
$ws.Cells.Item(23, 89).Value = 1
$ws.Cells.Item(23, 90).Value = 0.000165924409022183
$ws.Cells.Item(23, 91).Value = 0.0375737567626968

$ws.Cells.Item(24, 86).Value = 1065.30366649603
$ws.Cells.Item(24, 87).Value = 0.030822161267925
$ws.Cells.Item(24, 88).Value = 0.239863882163611
$ws.Cells.Item(24, 89).Value = 1
$ws.Cells.Item(24, 90).Value = 0.00754703573865318
$ws.Cells.Item(24, 91).Value = 0.114872293620683

$ws.Cells.Item(25, 86).Value = 897.65266997753
$ws.Cells.Item(25, 87).Value = 0.0228419717683309
$ws.Cells.Item(25, 88).Value = 0.250774243983707
$ws.Cells.Item(25, 89).Value = 1
$ws.Cells.Item(25, 90).Value = 0.00461774164710659
$ws.Cells.Item(25, 91).Value = 0.100181559443678

$ws.Cells.Item(26, 86).Value = 1039.98845870869
$ws.Cells.Item(26, 87).Value = 0.0252426327053718
$ws.Cells.Item(26, 88).Value = 0.252628224832543
$ws.Cells.Item(26, 89).Value = 1
$ws.Cells.Item(26, 90).Value = 0.00547264855283658
$ws.Cells.Item(26, 91).Value = 0.0387934060939835

$ws.Cells.Item(27, 86).Value = 945.317848171356
$ws.Cells.Item(27, 87).Value = 0.0198735016474076
$ws.Cells.Item(27, 88).Value = 0.193460669067578
$ws.Cells.Item(27, 89).Value = 0.899083749185418
$ws.Cells.Item(27, 90).Value = 0.00260755269656311
$ws.Cells.Item(27, 91).Value = 0.837755833281827

$ws.Cells.Item(28, 86).Value = 750.554774106368
$ws.Cells.Item(28, 87).Value = 0.0145054898283571
$ws.Cells.Item(28, 88).Value = 0.470835215982707
$ws.Cells.Item(28, 89).Value = 1
$ws.Cells.Item(28, 90).Value = 0.0129698748648404
$ws.Cells.Item(28, 91).Value = 0.373546634421696

$ws.Cells.Item(29, 86).Value = 1813.06397875218
$ws.Cells.Item(29, 87).Value = 0.042279189467343
$ws.Cells.Item(29, 88).Value = 0.152995216602583
$ws.Cells.Item(29, 89).Value = 1
$ws.Cells.Item(29, 90).Value = 0.00237059411817699
$ws.Cells.Item(29, 91).Value = 0.977446443129766

$ws.Cells.Item(30, 86).Value = 693.157549568385
$ws.Cells.Item(30, 87).Value = 0.0209251486238053
$ws.Cells.Item(30, 88).Value = 0.381273041581858
$ws.Cells.Item(30, 89).Value = 1
$ws.Cells.Item(30, 90).Value = 0.00718852001868104
$ws.Cells.Item(30, 91).Value = 0.756859315112003

$ws.Cells.Item(31, 86).Value = 780.737594022907
$ws.Cells.Item(31, 87).Value = 0.0230255987583002
$ws.Cells.Item(31, 88).Value = 0.318167292093443
$ws.Cells.Item(31, 89).Value = 1
$ws.Cells.Item(31, 90).Value = 0.00356605919047387
$ws.Cells.Item(31, 91).Value = 0.235825878617393

$ws.Cells.Item(32, 86).Value = 1255.57621289841
$ws.Cells.Item(32, 87).Value = 0.0411620072266793
$ws.Cells.Item(32, 88).Value = 0.304421589419485
$ws.Cells.Item(32, 89).Value = 1
$ws.Cells.Item(32, 90).Value = 0.00376293218012888
$ws.Cells.Item(32, 91).Value = 0.602252407504304

$ws.Cells.Item(33, 86).Value = 1337.54840412283
$ws.Cells.Item(33, 87).Value = 0.0783533340331162
$ws.Cells.Item(33, 88).Value = 0.270002944995831
$ws.Cells.Item(33, 89).Value = 1
$ws.Cells.Item(33, 90).Value = 0.00309094622639685
$ws.Cells.Item(33, 91).Value = 0.0142872321666929

$ws.Cells.Item(34, 86).Value = 900.816644188244
$ws.Cells.Item(34, 87).Value = 0.0371407981634906
$ws.Cells.Item(34, 88).Value = 0.323386235929562
$ws.Cells.Item(34, 89).Value = 1
$ws.Cells.Item(34, 90).Value = 0.00847977429642553
$ws.Cells.Item(34, 91).Value = 0.0685700553573564

$ws.Cells.Item(35, 86).Value = 1586.62756881026
$ws.Cells.Item(35, 87).Value = 0.00302738752844989
$ws.Cells.Item(35, 88).Value = 79.2623779795883
$ws.Cells.Item(35, 89).Value = 0
$ws.Cells.Item(35, 90).Value = 0
$ws.Cells.Item(35, 91).Value = 0.0000449734466792929

$ws.Cells.Item(36, 86).Value = 670.159780691187
$ws.Cells.Item(36, 87).Value = 0.00887362306503985
$ws.Cells.Item(36, 88).Value = 0.457360499071905
$ws.Cells.Item(36, 89).Value = 0.932770876130691
$ws.Cells.Item(36, 90).Value = 0.0110694019925257
$ws.Cells.Item(36, 91).Value = 1

$ws.Cells.Item(37, 86).Value = 1171.7127333449
$ws.Cells.Item(37, 87).Value = 0.027176396238687
$ws.Cells.Item(37, 88).Value = 0.276560658856001
$ws.Cells.Item(37, 89).Value = 1
$ws.Cells.Item(37, 90).Value = 0.00333870605169812
$ws.Cells.Item(37, 91).Value = 0.0412386496356317

$ws.Cells.Item(38, 86).Value = 744.757731221336
$ws.Cells.Item(38, 87).Value = 0.0428550249783385
$ws.Cells.Item(38, 88).Value = 0.244214602264428
$ws.Cells.Item(38, 89).Value = 1
$ws.Cells.Item(38, 90).Value = 0.00903360284465825
$ws.Cells.Item(38, 91).Value = 0.647279420404652

# --- Drop the forget_attention parameter: clear out column CN (header + data) ---
$ws.Range("CN2:CN38").ClearContents()

# --- AICc now penalizes for k=5 free parameters instead of k=6 ---
$ws.Range("CO3:CO38").Formula = "=CH3+2*5*(144/(144-5-1))"

# --- CL4 no longer needs the scientific-notation number format (was s="1") ---
$ws.Range("CL4").ClearFormats()

# --- Column widths auto-fit after the forget_attention column's data was cleared ---
$ws.Columns.Item(91).AutoFit()
$ws.Columns.Item(92).AutoFit()

# --- Viewport: scroll right and move the selection, as left by the author ---
$excel.ActiveWindow.ScrollColumn = 22
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AE20").Select()
